$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.459557
$ws.Range("H2").Value = 25.378671
$ws.Range("I2").Value = 0.3030961495696597
$ws.Range("J2").Value = 0.3030961495696597
$ws.Range("M2").Value = 8.839952666666667
$ws.Range("N2").Value = 26.519858
$ws.Range("O2").Value = 0.1312786621044678
$ws.Range("P2").Value = 0.1312786621044678
$ws.Range("Q2").Value = 74.78208346096868
$ws.Range("R2").Value = 673.038751148718
$ws.Range("S2").Value = 0.0397900570045206
$ws.Range("T2").Value = 0.03979005700452059
$ws.Range("G3").Value = 8.459557
$ws.Range("H3").Value = 25.378671
$ws.Range("I3").Value = 0.3030961495696597
$ws.Range("J3").Value = 0.3030961495696597
$ws.Range("O3").Value = 0.04535489698218144
$ws.Range("P3").Value = 0.04535489698218144
$ws.Range("Q3").Value = 25.836138463889
$ws.Range("R3").Value = 232.525246175001
$ws.Range("S3").Value = 0.01374689463942778
$ws.Range("T3").Value = 0.01374689463942777
$ws.Range("G4").Value = 8.459557
$ws.Range("H4").Value = 25.378671
$ws.Range("I4").Value = 0.3030961495696597
$ws.Range("J4").Value = 0.3030961495696597
$ws.Range("M4").Value = 13.31611333333333
$ws.Range("N4").Value = 39.94834
$ws.Range("O4").Value = 0.19775236460521
$ws.Range("P4").Value = 0.19775236460521
$ws.Range("Q4").Value = 112.6484197617933
$ws.Range("R4").Value = 1013.83577785614
$ws.Range("S4").Value = 0.05993798028013463
$ws.Range("T4").Value = 0.05993798028013462
$ws.Range("G5").Value = 8.459557
$ws.Range("H5").Value = 25.378671
$ws.Range("I5").Value = 0.3030961495696597
$ws.Range("J5").Value = 0.3030961495696597
$ws.Range("M5").Value = 30.323401
$ws.Range("N5").Value = 90.970203
$ws.Range("O5").Value = 0.4503209082496537
$ws.Range("P5").Value = 0.4503209082496537
$ws.Range("Q5").Value = 256.522539193357
$ws.Range("R5").Value = 2308.702852740213
$ws.Range("S5").Value = 0.1364905333611821
$ws.Range("T5").Value = 0.136490533361182
$ws.Range("G6").Value = 8.459557
$ws.Range("H6").Value = 25.378671
$ws.Range("I6").Value = 0.3030961495696597
$ws.Range("J6").Value = 0.3030961495696597
$ws.Range("M6").Value = 11.80377133333333
$ws.Range("N6").Value = 35.411314
$ws.Range("O6").Value = 0.175293168058487
$ws.Range("P6").Value = 0.175293168058487
$ws.Range("Q6").Value = 99.85467640929933
$ws.Range("R6").Value = 898.6920876836939
$ws.Range("S6").Value = 0.05313068428439467
$ws.Range("T6").Value = 0.05313068428439467
$ws.Range("I7").Value = 0.4601547065605718
$ws.Range("J7").Value = 0.4601547065605718
$ws.Range("M7").Value = 8.839952666666667
$ws.Range("N7").Value = 26.519858
$ws.Range("O7").Value = 0.1312786621044678
$ws.Range("P7").Value = 0.1312786621044678
$ws.Range("Q7").Value = 113.5327113849118
$ws.Range("R7").Value = 1021.794402464206
$ws.Range("S7").Value = 0.06040849423834585
$ws.Range("T7").Value = 0.06040849423834584
$ws.Range("I8").Value = 0.4601547065605718
$ws.Range("J8").Value = 0.4601547065605718
$ws.Range("O8").Value = 0.04535489698218144
$ws.Range("P8").Value = 0.04535489698218144
$ws.Range("S8").Value = 0.02087026931192067
$ws.Range("T8").Value = 0.02087026931192066
$ws.Range("I9").Value = 0.4601547065605718
$ws.Range("J9").Value = 0.4601547065605718
$ws.Range("M9").Value = 13.31611333333333
$ws.Range("N9").Value = 39.94834
$ws.Range("O9").Value = 0.19775236460521
$ws.Range("P9").Value = 0.19775236460521
$ws.Range("Q9").Value = 171.0206500927089
$ws.Range("R9").Value = 1539.18585083438
$ws.Range("S9").Value = 0.09099668130656964
$ws.Range("T9").Value = 0.09099668130656963
$ws.Range("I10").Value = 0.4601547065605718
$ws.Range("J10").Value = 0.4601547065605718
$ws.Range("M10").Value = 30.323401
$ws.Range("N10").Value = 90.970203
$ws.Range("O10").Value = 0.4503209082496537
$ws.Range("P10").Value = 0.4503209082496537
$ws.Range("Q10").Value = 389.4475529177357
$ws.Range("R10").Value = 3505.027976259621
$ws.Range("S10").Value = 0.2072172853937096
$ws.Range("T10").Value = 0.2072172853937096
$ws.Range("I11").Value = 0.4601547065605718
$ws.Range("J11").Value = 0.4601547065605718
$ws.Range("M11").Value = 11.80377133333333
$ws.Range("N11").Value = 35.411314
$ws.Range("O11").Value = 0.175293168058487
$ws.Range("P11").Value = 0.175293168058487
$ws.Range("Q11").Value = 151.5974366123109
$ws.Range("R11").Value = 1364.376929510798
$ws.Range("S11").Value = 0.08066197631002607
$ws.Range("T11").Value = 0.08066197631002607
$ws.Range("G12").Value = 1.955432333333333
$ws.Range("H12").Value = 5.866296999999999
$ws.Range("I12").Value = 0.0700608803720276
$ws.Range("J12").Value = 0.0700608803720276
$ws.Range("M12").Value = 8.839952666666667
$ws.Range("N12").Value = 26.519858
$ws.Range("O12").Value = 0.1312786621044678
$ws.Range("P12").Value = 0.1312786621044678
$ws.Range("Q12").Value = 17.28592926953622
$ws.Range("R12").Value = 155.573363425826
$ws.Range("S12").Value = 0.009197498641100952
$ws.Range("T12").Value = 0.009197498641100952
$ws.Range("G13").Value = 1.955432333333333
$ws.Range("H13").Value = 5.866296999999999
$ws.Range("I13").Value = 0.0700608803720276
$ws.Range("J13").Value = 0.0700608803720276
$ws.Range("O13").Value = 0.04535489698218144
$ws.Range("P13").Value = 0.04535489698218144
$ws.Range("Q13").Value = 5.972040914289666
$ws.Range("R13").Value = 53.748368228607
$ws.Range("S13").Value = 0.003177604011754249
$ws.Range("T13").Value = 0.003177604011754249
$ws.Range("G14").Value = 1.955432333333333
$ws.Range("H14").Value = 5.866296999999999
$ws.Range("I14").Value = 0.0700608803720276
$ws.Range("J14").Value = 0.0700608803720276
$ws.Range("M14").Value = 13.31611333333333
$ws.Range("N14").Value = 39.94834
$ws.Range("O14").Value = 0.19775236460521
$ws.Range("P14").Value = 0.19775236460521
$ws.Range("Q14").Value = 26.03875856633111
$ws.Range("R14").Value = 234.34882709698
$ws.Range("S14").Value = 0.0138547047598912
$ws.Range("T14").Value = 0.0138547047598912
$ws.Range("G15").Value = 1.955432333333333
$ws.Range("H15").Value = 5.866296999999999
$ws.Range("I15").Value = 0.0700608803720276
$ws.Range("J15").Value = 0.0700608803720276
$ws.Range("M15").Value = 30.323401
$ws.Range("N15").Value = 90.970203
$ws.Range("O15").Value = 0.4503209082496537
$ws.Range("P15").Value = 0.4503209082496537
$ws.Range("Q15").Value = 59.29535877203233
$ws.Range("R15").Value = 533.6582289482909
$ws.Range("S15").Value = 0.0315498792819018
$ws.Range("T15").Value = 0.0315498792819018
$ws.Range("G16").Value = 1.955432333333333
$ws.Range("H16").Value = 5.866296999999999
$ws.Range("I16").Value = 0.0700608803720276
$ws.Range("J16").Value = 0.0700608803720276
$ws.Range("M16").Value = 11.80377133333333
$ws.Range("N16").Value = 35.411314
$ws.Range("O16").Value = 0.175293168058487
$ws.Range("P16").Value = 0.175293168058487
$ws.Range("Q16").Value = 23.08147612047311
$ws.Range("R16").Value = 207.733285084258
$ws.Range("S16").Value = 0.01228119367737938
$ws.Range("T16").Value = 0.01228119367737939
$ws.Range("G17").Value = 2.929608
$ws.Range("H17").Value = 8.788824000000002
$ws.Range("I17").Value = 0.1049644685352285
$ws.Range("J17").Value = 0.1049644685352285
$ws.Range("M17").Value = 8.839952666666667
$ws.Range("N17").Value = 26.519858
$ws.Range("O17").Value = 0.1312786621044678
$ws.Range("P17").Value = 0.1312786621044678
$ws.Range("Q17").Value = 25.89759605188801
$ws.Range("R17").Value = 233.078364466992
$ws.Range("S17").Value = 0.0137795949978113
$ws.Range("T17").Value = 0.0137795949978113
$ws.Range("G18").Value = 2.929608
$ws.Range("H18").Value = 8.788824000000002
$ws.Range("I18").Value = 0.1049644685352285
$ws.Range("J18").Value = 0.1049644685352285
$ws.Range("O18").Value = 0.04535489698218144
$ws.Range("P18").Value = 0.04535489698218144
$ws.Range("Q18").Value = 8.947248411816002
$ws.Range("R18").Value = 80.52523570634402
$ws.Range("S18").Value = 0.004760652657204712
$ws.Range("T18").Value = 0.004760652657204712
$ws.Range("G19").Value = 2.929608
$ws.Range("H19").Value = 8.788824000000002
$ws.Range("I19").Value = 0.1049644685352285
$ws.Range("J19").Value = 0.1049644685352285
$ws.Range("M19").Value = 13.31611333333333
$ws.Range("N19").Value = 39.94834
$ws.Range("O19").Value = 0.19775236460521
$ws.Range("P19").Value = 0.19775236460521
$ws.Range("Q19").Value = 39.01099215024001
$ws.Range("R19").Value = 351.0989293521601
$ws.Range("S19").Value = 0.0207569718523706
$ws.Range("T19").Value = 0.0207569718523706
$ws.Range("G20").Value = 2.929608
$ws.Range("H20").Value = 8.788824000000002
$ws.Range("I20").Value = 0.1049644685352285
$ws.Range("J20").Value = 0.1049644685352285
$ws.Range("M20").Value = 30.323401
$ws.Range("N20").Value = 90.970203
$ws.Range("O20").Value = 0.4503209082496537
$ws.Range("P20").Value = 0.4503209082496537
$ws.Range("Q20").Value = 88.83567815680802
$ws.Range("R20").Value = 799.5211034112722
$ws.Range("S20").Value = 0.04726769480472628
$ws.Range("T20").Value = 0.04726769480472628
$ws.Range("G21").Value = 2.929608
$ws.Range("H21").Value = 8.788824000000002
$ws.Range("I21").Value = 0.1049644685352285
$ws.Range("J21").Value = 0.1049644685352285
$ws.Range("M21").Value = 11.80377133333333
$ws.Range("N21").Value = 35.411314
$ws.Range("O21").Value = 0.175293168058487
$ws.Range("P21").Value = 0.175293168058487
$ws.Range("Q21").Value = 34.580422928304
$ws.Range("R21").Value = 311.223806354736
$ws.Range("S21").Value = 0.01839955422311557
$ws.Range("T21").Value = 0.01839955422311558
$ws.Range("G22").Value = 1.722740333333333
$ws.Range("H22").Value = 5.168221
$ws.Range("I22").Value = 0.06172379496251228
$ws.Range("J22").Value = 0.06172379496251227
$ws.Range("M22").Value = 8.839952666666667
$ws.Range("N22").Value = 26.519858
$ws.Range("O22").Value = 0.1312786621044678
$ws.Range("P22").Value = 0.1312786621044678
$ws.Range("Q22").Value = 15.22894300362422
$ws.Range("R22").Value = 137.060487032618
$ws.Range("S22").Value = 0.008103017222689102
$ws.Range("T22").Value = 0.0081030172226891
$ws.Range("G23").Value = 1.722740333333333
$ws.Range("H23").Value = 5.168221
$ws.Range("I23").Value = 0.06172379496251228
$ws.Range("J23").Value = 0.06172379496251227
$ws.Range("O23").Value = 0.04535489698218144
$ws.Range("P23").Value = 0.04535489698218144
$ws.Range("Q23").Value = 5.261381629005666
$ws.Range("R23").Value = 47.352434661051
$ws.Range("S23").Value = 0.002799476361874034
$ws.Range("T23").Value = 0.002799476361874034
$ws.Range("G24").Value = 1.722740333333333
$ws.Range("H24").Value = 5.168221
$ws.Range("I24").Value = 0.06172379496251228
$ws.Range("J24").Value = 0.06172379496251227
$ws.Range("M24").Value = 13.31611333333333
$ws.Range("N24").Value = 39.94834
$ws.Range("O24").Value = 0.19775236460521
$ws.Range("P24").Value = 0.19775236460521
$ws.Range("Q24").Value = 22.94020552257111
$ws.Range("R24").Value = 206.46184970314
$ws.Range("S24").Value = 0.01220602640624395
$ws.Range("T24").Value = 0.01220602640624395
$ws.Range("G25").Value = 1.722740333333333
$ws.Range("H25").Value = 5.168221
$ws.Range("I25").Value = 0.06172379496251228
$ws.Range("J25").Value = 0.06172379496251227
$ws.Range("M25").Value = 30.323401
$ws.Range("N25").Value = 90.970203
$ws.Range("O25").Value = 0.4503209082496537
$ws.Range("P25").Value = 0.4503209082496537
$ws.Range("Q25").Value = 52.23934594654033
$ws.Range("R25").Value = 470.154113518863
$ws.Range("S25").Value = 0.02779551540813393
$ws.Range("T25").Value = 0.02779551540813393
$ws.Range("G26").Value = 1.722740333333333
$ws.Range("H26").Value = 5.168221
$ws.Range("I26").Value = 0.06172379496251228
$ws.Range("J26").Value = 0.06172379496251227
$ws.Range("M26").Value = 11.80377133333333
$ws.Range("N26").Value = 35.411314
$ws.Range("O26").Value = 0.175293168058487
$ws.Range("P26").Value = 0.175293168058487
$ws.Range("Q26").Value = 20.33483296137711
$ws.Range("R26").Value = 183.013496652394
$ws.Range("S26").Value = 0.01081975956357126
$ws.Range("T26").Value = 0.01081975956357126
